$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "Alterar o readme file no repositorio git" from F4 down to G8
$ws.Range("G8").Value = $ws.Range("F4").Value()
$ws.Range("F4").Value = $null

# Update the active selection to D4
$ws.Range("D4").Select()
